# Generate Report for Archive
#
# 1) Status text "Ready for handoff" -> "In Translation" on every sheet that
#    surfaces it (Overview!E2:F3, zh-cn!C2:C3, de-de!C2:C3).
# 2) The Status-ish columns (Overview E:F, zh-cn C, de-de C) are narrowed to
#    match the now-shorter text ("In Translation" is shorter than
#    "Ready for handoff"), mirroring an Excel "AutoFit columns" pass.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# --- Overview sheet -------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# Narrow columns E & F to fit the shorter text.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# --- zh-cn sheet ------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

# --- de-de sheet --------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
